$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header row (bold, centered, bordered)
# by copying the formatting from AC1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the team record (Wins/Losses/Ties) for every player row.
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}

$wb.Save()
